# Actualización automática 2025-11-17 17:30:08
$wb = $excel.ActiveWorkbook

$wsGrupo    = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual  = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumpl    = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# ---- Sheet 1: VENTAS POR GRUPO ----
$wsGrupo.Range("D24").Value = 933.12
$wsGrupo.Range("L29").Value = 2215.01
$wsGrupo.Range("M29").Value = 617.41
$wsGrupo.Range("I31").Value = 79.2
$wsGrupo.Range("L48").Value = 886.88
$wsGrupo.Range("M52").Value = 616.84
$wsGrupo.Range("L56").Value = "7 de 54"
$wsGrupo.Range("M56").Value = "13 de 54"

# ---- Sheet 2: VENTA MENSUAL ----
$wsMensual.Range("F24").Value = 4381.16
$wsMensual.Range("F29").Value = 4432.98
$wsMensual.Range("F31").Value = 79.2
$wsMensual.Range("F48").Value = 1362.08
$wsMensual.Range("F53").Value = 616.84
$wsMensual.Range("F54").Value = 616.84
$wsMensual.Range("F60").Value = 28133.16

# ---- Sheet 3: CUMPLIMIENTO MENSUAL ----
$wsCumpl.Range("D3").Value = 3818.98
$wsCumpl.Range("E3").Value = 2804.28
$wsCumpl.Range("F3").Value = 0.5766012507435915

$wsCumpl.Range("D7").Value = 415.8
$wsCumpl.Range("E7").Value = 904.2
$wsCumpl.Range("F7").Value = 0.315

$wsCumpl.Range("D11").Value = 10166.24
$wsCumpl.Range("E11").Value = 4069.75
$wsCumpl.Range("F11").Value = 0.7141224459977845

$wsCumpl.Range("D12").Value = 10312.76
$wsCumpl.Range("E12").Value = 54631.24
$wsCumpl.Range("F12").Value = 0.1587946538556295

$wsCumpl.Range("D14").Value = 27492.92
$wsCumpl.Range("E14").Value = 71463.33685923838
$wsCumpl.Range("F14").Value = 0.2778290213534215

# Column width tweaks on CUMPLIMIENTO MENSUAL (D and E) matching auto-fit after value update
# (the runtime stores width as ColumnWidth + 5/6, so compensate to land on the exact target width)
$wsCumpl.Columns.Item(4).ColumnWidth = 13.166666666666666
$wsCumpl.Columns.Item(5).ColumnWidth = 21.166666666666668
